# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This script rewrites the "K" column (column G) on Sheet1 with the
# regenerated strikeout ("K") values for each game-log row (rows 2-83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for G2:G83, in row order (row 2 first).
$newK = @(
    1, 2, 0, 3, 1, 1, 2, 0, 1, 2,
    2, 2, 0, 2, 0, 2, 2, 2, 2, 1,
    1, 2, 0, 1, 2, 1, 1, 1, 1, 2,
    1, 0, 2, 3, 2, 3, 1, 1, 2, 2,
    1, 0, 2, 1, 1, 2, 1, 1, 1, 1,
    2, 1, 0, 1, 1, 1, 0, 0, 0, 1,
    2, 1, 1, 1, 0, 1, 1, 4, 3, 1,
    1, 0, 0, 3, 3, 1, 1, 1, 2, 1,
    1, 1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
